$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 16:29"

# --- Reorder country names (content of the shared strings used by these rows) ---
# Costa Rica / Portugal swap position (rows 52-53)
$ws.Range("A52").Value = "Portugal"
$ws.Range("A53").Value = "Costa Rica"

# Birmania moves ahead of Noruega/Albania/Sudan (rows 93-96)
$ws.Range("A93").Value = "Birmania"
$ws.Range("A94").Value = "Noruega"
$ws.Range("A95").Value = "Albania"
$ws.Range("A96").Value = "Sudan"

# Santa Lucia / Nueva Caledonia swap position (rows 207-208)
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Update daily statistics for affected rows (columns B,C,D,E,F,G,H) ---
# Row 4
$ws.Range("B4").Value = 7454758
$ws.Range("C4").Value = 7476
$ws.Range("D4").Value = 4701444
$ws.Range("E4").Value = 2541465
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 109
$ws.Range("H4").Value = 211849

# Row 5
$ws.Range("B5").Value = 6323247
$ws.Range("C5").Value = 12980
$ws.Range("D5").Value = 5280204
$ws.Range("E5").Value = 944221
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 114
$ws.Range("H5").Value = 98822

# Row 52
$ws.Range("B52").Value = 76396
$ws.Range("C52").Value = 854
$ws.Range("D52").Value = 48937
$ws.Range("E52").Value = 25482
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 6
$ws.Range("H52").Value = 1977

# Row 53
$ws.Range("B53").Value = 75760
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 37841
$ws.Range("E53").Value = 37015
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 904

# Row 61
$ws.Range("B61").Value = 53832
$ws.Range("C61").Value = 550
$ws.Range("D61").Value = 45300
$ws.Range("E61").Value = 6458
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 2074

# Row 66
$ws.Range("B66").Value = 46656
$ws.Range("C66").Value = 30
$ws.Range("D66").Value = 45942
$ws.Range("E66").Value = 413
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 301

# Row 73
$ws.Range("B73").Value = 38713
$ws.Range("C73").Value = 184
$ws.Range("D73").Value = 24908
$ws.Range("E73").Value = 13094
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 711

# Row 79
$ws.Range("B79").Value = 27749
$ws.Range("C79").Value = 280
$ws.Range("D79").Value = 20947
$ws.Range("E79").Value = 5941
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 5
$ws.Range("H79").Value = 861

# Row 88
$ws.Range("B88").Value = 18138
$ws.Range("C88").Value = 161
$ws.Range("D88").Value = 15068
$ws.Range("E88").Value = 2327
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 743

# Row 93
$ws.Range("B93").Value = 14383
$ws.Range("C93").Value = 1010
$ws.Range("D93").Value = 4156
$ws.Range("E93").Value = 9906
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 11
$ws.Range("H93").Value = 321

# Row 94
$ws.Range("B94").Value = 14085
$ws.Range("C94").Value = 58
$ws.Range("D94").Value = 11190
$ws.Range("E94").Value = 2621
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 274

# Row 95
$ws.Range("B95").Value = 13649
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 7847
$ws.Range("E95").Value = 5415
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 387

# Row 96
$ws.Range("B96").Value = 13640
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 6764
$ws.Range("E96").Value = 6040
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 836

# Row 107
$ws.Range("B107").Value = 9811
$ws.Range("C107").Value = 42
$ws.Range("D107").Value = 8617
$ws.Range("E107").Value = 1117
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 77

# Row 120
$ws.Range("B120").Value = 5670
$ws.Range("C120").Value = 73
$ws.Range("D120").Value = 4922
$ws.Range("E120").Value = 626
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 122

# Row 143
$ws.Range("B143").Value = 3380
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 3233
$ws.Range("E143").Value = 134
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 13

# Row 174
$ws.Range("B174").Value = 539
$ws.Range("C174").Value = 5
$ws.Range("D174").Value = 519
$ws.Range("E174").Value = 13
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 7

